$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, copying the formatting used by the other
# header cells (bold, bordered, centered) from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new column: 0 for every row except row 12 which is 1
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
